$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column) - this shifts the
# existing "Late"/"heading"/"Outstanding" columns one place to the right,
# matching the new "Variable Instalments" layout.
$ws.Columns("N").Insert()

# The newly inserted column should carry the same width as its neighbours
# (11 characters, as stored in the sheet's <cols> definition).
$ws.Columns("N").ColumnWidth = 10.140625

# Make "Repayment schedule" the active sheet/tab.
$ws.Activate()

# Update the current selection on the Repayment schedule sheet.
$ws.Range("R7").Select()
